# Refresh the cryptocurrency price/volume snapshot on Sheet1, as produced by
# the scheduled GitHub Actions scraper run on Thu Feb 16 12:51:35 UTC 2023.
#
# All data cells in this sheet are stored as literal text (e.g. "319.39",
# "5.71%") rather than numbers/percentages, so each write forces the Text
# number format before assigning the value (otherwise Excel would silently
# reinterpret "319.39" as a float or "5.71%" as a percentage fraction), then
# restores the cell to the workbook's default "Normal" style so no stray
# formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-Text($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2 - BNB
Set-Text "D2" "319.39"
Set-Text "E2" "5.71%"

# Row 3 - OKB
Set-Text "D3" "48.53"
Set-Text "E3" "13.97%"

# Row 4 - HuobiToken
Set-Text "D4" "5.248"
Set-Text "E4" "4.53%"

# Row 5 - Cronos
Set-Text "D5" "0.08074"
Set-Text "E5" "5.10%"

# Row 6 - GateToken
Set-Text "D6" "4.572"
Set-Text "E6" "4.34%"

# Row 7 - was FTXToken, now MXToken (rows 7 & 8 swapped identity)
Set-Text "B7" "MXToken"
Set-Text "C7" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-Text "D7" "1.249"
Set-Text "E7" "25.24%"

# Row 8 - was MXToken, now FTXToken
Set-Text "B8" "FTXToken"
Set-Text "C8" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-Text "D8" "1.650"
Set-Text "E8" "2.26%"

# Row 9 - LiechtensteinCryptoassetsExchange
Set-Text "D9" "0.1295"
Set-Text "E9" "7.20%"

# Row 10 - WazirX
Set-Text "D10" "0.1927"
Set-Text "E10" "4.17%"

# Row 11 - MandalaExchangeToken
Set-Text "D11" "0.09477"
Set-Text "E11" "4.75%"

# Row 12 - BitrueCoin
Set-Text "D12" "0.04592"
Set-Text "E12" "12.26%"

# Row 13 - BitMartToken
Set-Text "D13" "0.1044"
Set-Text "E13" "0.10%"

# Row 14 - BitForexToken
Set-Text "D14" "0.001335"
Set-Text "E14" "5.18%"

# Row 15 - CoinExToken
Set-Text "D15" "0.04158"
Set-Text "E15" "0.36%"

# Row 16 - TigerCash
Set-Text "D16" "0.005805"
Set-Text "E16" "-1.96%"

# Row 17 - LEO
Set-Text "D17" "3.340"
Set-Text "E17" "-0.63%"

# Row 18 - BTSEToken
Set-Text "D18" "2.437"
Set-Text "E18" "2.22%"

# Row 19 - BitpandaEcosystemToken
Set-Text "D19" "0.3405"
Set-Text "E19" "1.94%"

# Row 20 - MCDex
Set-Text "D20" "8.205"
Set-Text "E20" "-1.87%"

# Row 21 - ProBitToken
Set-Text "D21" "0.1390"
Set-Text "E21" "2.15%"

# Row 22 - ZBToken
Set-Text "D22" "0.3097"
Set-Text "E22" "3.76%"

# Row 23 - BitKan
Set-Text "D23" "0.001306"
Set-Text "E23" "3.35%"

# Row 24 - HotbitToken
Set-Text "D24" "0.004249"
Set-Text "E24" "7.38%"

# Row 25 - NitroEx
Set-Text "D25" "0.0001350"
Set-Text "E25" "0.49%"

# Row 26 - UpBots
Set-Text "D26" "0.0003538"
Set-Text "E26" "-95.24%"

# Row 38
Set-Text "D38" "0.02705"
Set-Text "E38" "9.76%"

# Row 39 - IDEX
Set-Text "D39" "0.05678"
Set-Text "E39" "7.56%"

# Row 40 - CEJI (only volume% changed)
Set-Text "E40" "2.13%"

# Row 41 - KickToken
Set-Text "D41" "0.007987"
Set-Text "E41" "4.59%"

# Row 42 - BKEXToken
Set-Text "D42" "0.1442"
Set-Text "E42" "7.10%"

# Row 43 - Dexo
Set-Text "D43" "0.007685"
Set-Text "E43" "4.68%"

# Row 44 - LocalTraders
Set-Text "D44" "0.008707"
Set-Text "E44" "18.97%"

# Row 45 - PooCoin
Set-Text "D45" "0.3510"
Set-Text "E45" "15.54%"

# Row 46 - CoinLion
Set-Text "D46" "0.00006865"
Set-Text "E46" "4.68%"

# Row 47 - Kangarootoken
Set-Text "D47" "0.00000000750"
Set-Text "E47" "0.67%"

# Row 48 - BOLO
Set-Text "D48" "0.05436"
Set-Text "E48" "18.33%"

# Row 49 - CoinbaseStockToken
Set-Text "D49" "0.003999"
Set-Text "E49" "-4.80%"

# Row 50 - CryptobidCoin
Set-Text "D50" "0.00002099"
Set-Text "E50" "0.67%"

# Row 51 - SpecialPowerGold
Set-Text "D51" "0.0002000"
Set-Text "E51" "0.67%"
